# Updated cryptos list on Sat Nov 23 15:27:54 UTC 2024 with GitHub Actions
#
# Applies the refreshed Price / Volume(1h) snapshot to the "cryptos" sheet:
#   - existing coins get new Price (col D) / Volume(1h) (col E) readings
#   - the coin list itself shifted (WrappedeETH dropped off, OKB newly
#     appeared at the bottom), so rows 29-51 get new Coin/Link/Price/Volume
#
# All Price/Volume cells in this sheet are stored as plain text (not
# numbers), even when their content looks numeric (e.g. "262.28" or
# "98.771.38" which isn't valid as a single number anyway because of the
# thousands-dot). Force NumberFormat to Text before writing so COM doesn't
# auto-convert number-looking strings to real numbers, then restore the
# "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '98.771.38' }
    @{ Cell = 'E2'; Value = '  +0.43%  ' }
    @{ Cell = 'D3'; Value = '3.469.86' }
    @{ Cell = 'E3'; Value = '  +5.08%  ' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '262.28' }
    @{ Cell = 'E5'; Value = '  +2.63%  ' }
    @{ Cell = 'D6'; Value = '675.48' }
    @{ Cell = 'E6'; Value = '  +8.62%  ' }
    @{ Cell = 'D7'; Value = '1.57' }
    @{ Cell = 'E7'; Value = '  +8.86%  ' }
    @{ Cell = 'D8'; Value = '0.461' }
    @{ Cell = 'E8'; Value = '  +14.69%  ' }
    @{ Cell = 'D9'; Value = '1.13' }
    @{ Cell = 'E9'; Value = '  +23.79%  ' }
    @{ Cell = 'E10'; Value = '  -0.05%  ' }
    @{ Cell = 'D11'; Value = '3.467.81' }
    @{ Cell = 'E11'; Value = '  +5.13%  ' }
    @{ Cell = 'D12'; Value = '0.219' }
    @{ Cell = 'E12'; Value = '  +10.11%  ' }
    @{ Cell = 'D13'; Value = '43.06' }
    @{ Cell = 'E13'; Value = '  +11.26%  ' }
    @{ Cell = 'E14'; Value = '  +10.43%  ' }
    @{ Cell = 'D15'; Value = '6.27' }
    @{ Cell = 'E15'; Value = '  +14.73%  ' }
    @{ Cell = 'D16'; Value = '98.331.29' }
    @{ Cell = 'E16'; Value = '  +0.34%  ' }
    @{ Cell = 'D17'; Value = '4.112.19' }
    @{ Cell = 'E17'; Value = '  +4.84%  ' }
    @{ Cell = 'D18'; Value = '8.74' }
    @{ Cell = 'E18'; Value = '  +41.06%  ' }
    @{ Cell = 'D19'; Value = '3.467.46' }
    @{ Cell = 'E19'; Value = '  +5.25%  ' }
    @{ Cell = 'D20'; Value = '17.75' }
    @{ Cell = 'E20'; Value = '  +17.41%  ' }
    @{ Cell = 'D21'; Value = '3.63' }
    @{ Cell = 'E21'; Value = '  +4.11%  ' }
    @{ Cell = 'D22'; Value = '530.45' }
    @{ Cell = 'E22'; Value = '  +10.11%  ' }
    @{ Cell = 'D23'; Value = '10.94' }
    @{ Cell = 'E23'; Value = '  +16.54%  ' }
    @{ Cell = 'D24'; Value = '0.492' }
    @{ Cell = 'E24'; Value = '  +68.25%  ' }
    @{ Cell = 'D25'; Value = '0.0000219' }
    @{ Cell = 'E25'; Value = '  +7.43%  ' }
    @{ Cell = 'D26'; Value = '6.49' }
    @{ Cell = 'E26'; Value = '  +16.47%  ' }
    @{ Cell = 'D27'; Value = '104.17' }
    @{ Cell = 'E27'; Value = '  +17.55%  ' }
    @{ Cell = 'D28'; Value = '13.09' }
    @{ Cell = 'E28'; Value = '  +10.70%  ' }
    @{ Cell = 'B29'; Value = 'Hedera' }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D29'; Value = '0.153' }
    @{ Cell = 'E29'; Value = '  +16.79%  ' }
    @{ Cell = 'B30'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D30'; Value = '11.70' }
    @{ Cell = 'E30'; Value = '  +17.33%  ' }
    @{ Cell = 'B31'; Value = 'Cronos' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D31'; Value = '0.198' }
    @{ Cell = 'E31'; Value = '  +5.62%  ' }
    @{ Cell = 'B32'; Value = 'Dai' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D32'; Value = '0.998' }
    @{ Cell = 'E32'; Value = '  -0.04%  ' }
    @{ Cell = 'B33'; Value = 'PolygonEcosystemToken' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol' }
    @{ Cell = 'D33'; Value = '0.584' }
    @{ Cell = 'E33'; Value = '  +27.97%  ' }
    @{ Cell = 'B34'; Value = 'EthereumClassic' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D34'; Value = '30.68' }
    @{ Cell = 'E34'; Value = '  +11.07%  ' }
    @{ Cell = 'D35'; Value = '0.999' }
    @{ Cell = 'E35'; Value = '  -0.03%  ' }
    @{ Cell = 'B36'; Value = 'PancakeSwap' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D36'; Value = '2.25' }
    @{ Cell = 'E36'; Value = '  +16.30%  ' }
    @{ Cell = 'B37'; Value = 'RenderToken' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render' }
    @{ Cell = 'D37'; Value = '8.08' }
    @{ Cell = 'E37'; Value = '  +12.49%  ' }
    @{ Cell = 'B38'; Value = 'Kaspa' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = 'D38'; Value = '0.162' }
    @{ Cell = 'E38'; Value = '  +9.64%  ' }
    @{ Cell = 'B39'; Value = 'Bittensor' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Cell = 'D39'; Value = '539.34' }
    @{ Cell = 'E39'; Value = '  +10.42%  ' }
    @{ Cell = 'B40'; Value = 'Fetch.AI' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D40'; Value = '1.43' }
    @{ Cell = 'E40'; Value = '  +16.05%  ' }
    @{ Cell = 'B41'; Value = 'WhiteBITCoin' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt' }
    @{ Cell = 'D41'; Value = '24.79' }
    @{ Cell = 'E41'; Value = '  -0.09%  ' }
    @{ Cell = 'B42'; Value = 'VeChain' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D42'; Value = '0.0445' }
    @{ Cell = 'E42'; Value = '  +37.54%  ' }
    @{ Cell = 'B43'; Value = 'ARBITRUM' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D43'; Value = '0.874' }
    @{ Cell = 'E43'; Value = '  +10.43%  ' }
    @{ Cell = 'B44'; Value = 'dogwifhat' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' }
    @{ Cell = 'D44'; Value = '3.51' }
    @{ Cell = 'E44'; Value = '  +12.72%  ' }
    @{ Cell = 'B45'; Value = 'MantraDAO' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om' }
    @{ Cell = 'D45'; Value = '3.74' }
    @{ Cell = 'E45'; Value = '  +2.50%  ' }
    @{ Cell = 'B46'; Value = 'Cosmos' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D46'; Value = '8.58' }
    @{ Cell = 'E46'; Value = '  +19.51%  ' }
    @{ Cell = 'B47'; Value = 'ImmutableX' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D47'; Value = '1.62' }
    @{ Cell = 'E47'; Value = '  +20.42%  ' }
    @{ Cell = 'B48'; Value = 'USDe' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' }
    @{ Cell = 'D48'; Value = '1.00' }
    @{ Cell = 'E48'; Value = '  +0.03%  ' }
    @{ Cell = 'B49'; Value = 'Filecoin' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D49'; Value = '5.34' }
    @{ Cell = 'E49'; Value = '  +15.91%  ' }
    @{ Cell = 'D50'; Value = '2.13' }
    @{ Cell = 'E50'; Value = '  +12.21%  ' }
    @{ Cell = 'B51'; Value = 'OKB' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D51'; Value = '52.06' }
    @{ Cell = 'E51'; Value = '  +14.74%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates to cryptos sheet"
